$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.882.84"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.892.37"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'0.7751"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'244.20"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.3143"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").Value = "'0.07408"
$ws.Range("E9").Value = "  +4.34%  "
$ws.Range("D10").Value = "'25.35"
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("D11").Value = "'0.08151"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "'0.7651"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "'5.474"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'92.53"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.791.01"
$ws.Range("E15").Value = "  -5.73%  "
$ws.Range("D16").Value = "'6.225"
$ws.Range("E16").Value = "  +5.07%  "
$ws.Range("D17").Value = "29.879.06"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'13.97"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'245.57"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "'0.000007871"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'8.142"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").Value = "2.126.32"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'0.1571"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").Value = "'9.438"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'162.25"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'2.042"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").Value = "'1.454"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("D31").Value = "'1.544"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'4.513"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "'0.05597"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "'4.097"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'1.251"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "'0.7594"
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("D37").Value = "'0.9995"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'2.646"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'2.790"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").Value = "1.147.93"
$ws.Range("E41").Value = "  +11.50%  "
$ws.Range("D42").Value = "'74.56"
$ws.Range("E42").Value = "  +2.94%  "
$ws.Range("D43").Value = "'0.4465"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'5.979"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").Value = "'0.8552"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'3.142"
$ws.Range("E48").Value = "  +6.20%  "
$ws.Range("D49").Value = "'101.98"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "'9.918"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'7.523"
$ws.Range("E51").Value = "  +0.33%  "

foreach ($addr in @("D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D16","D18","D19","D20","D21","D22","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D48","D49","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}
